$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.029649106843733
$ws.Range("D2").Value = 1.032851307206566
$ws.Range("E2").Value = 1.050504100679194
$ws.Range("F2").Value = 1.055443569404625
$ws.Range("I2").Value = 1.034135526793734
$ws.Range("J2").Value = 1.034794714404582
$ws.Range("K2").Value = 1.035655346612106
$ws.Range("L2").Value = 1.053258194707174
$ws.Range("M2").Value = 1.058184012321561
$ws.Range("N2").Value = 1.036264241591649

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.030454664576612
$ws.Range("D3").Value = 1.033427927279874
$ws.Range("E3").Value = 1.051712646146979
$ws.Range("F3").Value = 1.056708145356397
$ws.Range("I3").Value = 1.034286662105933
$ws.Range("J3").Value = 1.035242129110695
$ws.Range("K3").Value = 1.036041445900356
$ws.Range("L3").Value = 1.054278179386958
$ws.Range("M3").Value = 1.059260880229228
$ws.Range("N3").Value = 1.036712291677964

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.030976168715973
$ws.Range("D4").Value = 1.03380123363722
$ws.Range("E4").Value = 1.052495619388656
$ws.Range("F4").Value = 1.057527407178491
$ws.Range("I4").Value = 1.034383384410508
$ws.Range("J4").Value = 1.035531224320419
$ws.Range("K4").Value = 1.036290771536585
$ws.Range("L4").Value = 1.054938571842563
$ws.Range("M4").Value = 1.059958123164121
$ws.Range("N4").Value = 1.037001797436046

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.03119546872315
$ws.Range("D5").Value = 1.033958216791335
$ws.Range("E5").Value = 1.052825012646169
$ws.Range("F5").Value = 1.057872064057307
$ws.Range("I5").Value = 1.034423789237679
$ws.Range("J5").Value = 1.035652660446597
$ws.Range("K5").Value = 1.036395465886724
$ws.Range("L5").Value = 1.055216295707601
$ws.Range("M5").Value = 1.060251349241044
$ws.Range("N5").Value = 1.037123406015451

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.031232293622147
$ws.Range("D6").Value = 1.033984577560336
$ws.Range("E6").Value = 1.052880332806252
$ws.Range("F6").Value = 1.057929947509245
$ws.Range("I6").Value = 1.034430558287151
$ws.Range("J6").Value = 1.035673044246766
$ws.Range("K6").Value = 1.03641303733956
$ws.Range("L6").Value = 1.055262932345295
$ws.Range("M6").Value = 1.060300589368316
$ws.Range("N6").Value = 1.037143818762954

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.030979098780504
$ws.Range("D7").Value = 1.033803331078157
$ws.Range("E7").Value = 1.052500019847786
$ws.Range("F7").Value = 1.057532011559858
$ws.Range("I7").Value = 1.034383925313154
$ws.Range("J7").Value = 1.035532847347927
$ws.Range("K7").Value = 1.036292170948713
$ws.Range("L7").Value = 1.054942282428982
$ws.Range("M7").Value = 1.059962040851596
$ws.Range("N7").Value = 1.037003422768439

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.029921295425142
$ws.Range("D8").Value = 1.033046137506191
$ws.Range("E8").Value = 1.050912335160039
$ws.Range("F8").Value = 1.055870732779161
$ws.Range("I8").Value = 1.034186825411779
$ws.Range("J8").Value = 1.034946005093998
$ws.Range("K8").Value = 1.03578593498589
$ws.Range("L8").Value = 1.05360282199782
$ws.Range("M8").Value = 1.058547854567809
$ws.Range("N8").Value = 1.036415747131198

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.028059314083872
$ws.Range("D9").Value = 1.031713419171966
$ws.Range("E9").Value = 1.048121982580006
$ws.Range("F9").Value = 1.052950952678834
$ws.Range("I9").Value = 1.033831319697766
$ws.Range("J9").Value = 1.033908796392102
$ws.Range("K9").Value = 1.034890044028759
$ws.Range("L9").Value = 1.051245517422221
$ws.Range("M9").Value = 1.056059208019549
$ws.Range("N9").Value = 1.035377065473976

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.02681941303117
$ws.Range("D10").Value = 1.030826069319977
$ws.Range("E10").Value = 1.046266647465708
$ws.Range("F10").Value = 1.05100950513811
$ws.Range("I10").Value = 1.033588839000447
$ws.Range("J10").Value = 1.033215278305085
$ws.Range("K10").Value = 1.034290256165655
$ws.Range("L10").Value = 1.049675961139127
$ws.Range("M10").Value = 1.054402319944623
$ws.Range("N10").Value = 1.034682562511727

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.026282874032505
$ws.Range("D11").Value = 1.030442122557966
$ws.Range("E11").Value = 1.045464418972229
$ws.Range("F11").Value = 1.050170030199495
$ws.Range("I11").Value = 1.033482549351575
$ws.Range("J11").Value = 1.03291450257465
$ws.Range("K11").Value = 1.034029953122381
$ws.Range("L11").Value = 1.048996787349339
$ws.Range("M11").Value = 1.053685385917229
$ws.Range("N11").Value = 1.034381359645256

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.026083632770536
$ws.Range("D12").Value = 1.030299551358551
$ws.Range("E12").Value = 1.045166606485981
$ws.Range("F12").Value = 1.049858388871206
$ws.Range("I12").Value = 1.033442874723908
$ws.Range("J12").Value = 1.032802710100841
$ws.Range("K12").Value = 1.033933177074057
$ws.Range("L12").Value = 1.048744579444535
$ws.Range("M12").Value = 1.05341916051729
$ws.Range("N12").Value = 1.034269408413311

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.026126368270718
$ws.Range("D13").Value = 1.030330131363491
$ws.Range("E13").Value = 1.045230480557123
$ws.Range("F13").Value = 1.049925228997529
$ws.Range("I13").Value = 1.033451393833849
$ws.Range("J13").Value = 1.032826693164725
$ws.Range("K13").Value = 1.033953939828943
$ws.Range("L13").Value = 1.048798675855269
$ws.Range("M13").Value = 1.053476263349785
$ws.Range("N13").Value = 1.034293425535897

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.026266403595322
$ws.Range("D14").Value = 1.030430336678017
$ws.Range("E14").Value = 1.045439798218496
$ws.Range("F14").Value = 1.050144266223896
$ws.Range("I14").Value = 1.033479273787622
$ws.Range("J14").Value = 1.032905263217077
$ws.Range("K14").Value = 1.034021955372285
$ws.Range("L14").Value = 1.048975938393222
$ws.Range("M14").Value = 1.053663378100059
$ws.Range("N14").Value = 1.034372107166734

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.026352691066854
$ws.Range("D15").Value = 1.030492082315609
$ws.Range("E15").Value = 1.045568788354875
$ws.Range("F15").Value = 1.050279245719946
$ws.Range("I15").Value = 1.033496425865312
$ws.Range("J15").Value = 1.032953663438804
$ws.Range("K15").Value = 1.034063850370785
$ws.Range("L15").Value = 1.049085164611534
$ws.Range("M15").Value = 1.053778675698696
$ws.Range("N15").Value = 1.034420576122328

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.02685502875497
$ws.Range("D16").Value = 1.030851556655866
$ws.Range("E16").Value = 1.046319912746692
$ws.Range("F16").Value = 1.051065243204018
$ws.Range("I16").Value = 1.033595865873812
$ws.Range("J16").Value = 1.033235229801888
$ws.Range("K16").Value = 1.034307519245554
$ws.Range("L16").Value = 1.049721045173369
$ws.Range("M16").Value = 1.05444991116672
$ws.Range("N16").Value = 1.034702542341944

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.027170225471435
$ws.Range("D17").Value = 1.031077121802392
$ws.Range("E17").Value = 1.046791378798864
$ws.Range("F17").Value = 1.051558595160994
$ws.Range("I17").Value = 1.033657895883688
$ws.Range("J17").Value = 1.033411721589404
$ws.Range("K17").Value = 1.034460208699003
$ws.Range("L17").Value = 1.050120037354196
$ws.Range("M17").Value = 1.05487109532503
$ws.Range("N17").Value = 1.034879284768042

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.027354107721145
$ws.Range("D18").Value = 1.031208717209523
$ws.Range("E18").Value = 1.047066487399588
$ws.Range("F18").Value = 1.051846473375603
$ws.Range("I18").Value = 1.033693952063816
$ws.Range("J18").Value = 1.033514620105737
$ws.Range("K18").Value = 1.034549212798167
$ws.Range("L18").Value = 1.050352806642386
$ws.Range("M18").Value = 1.055116814120361
$ws.Range("N18").Value = 1.034982329412071

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.027416812435172
$ws.Range("D19").Value = 1.031253592391507
$ws.Range("E19").Value = 1.047160311082857
$ws.Range("F19").Value = 1.051944651834331
$ws.Range("I19").Value = 1.033706225094021
$ws.Range("J19").Value = 1.033549697986875
$ws.Range("K19").Value = 1.034579551188262
$ws.Range("L19").Value = 1.050432182424718
$ws.Range("M19").Value = 1.05520060624705
$ws.Range("N19").Value = 1.035017457107824

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.02713640438994
$ws.Range("D20").Value = 1.031052917986302
$ws.Range("E20").Value = 1.046740783537927
$ws.Range("F20").Value = 1.051505651342371
$ws.Range("I20").Value = 1.033651253563658
$ws.Range("J20").Value = 1.033392790459287
$ws.Range("K20").Value = 1.034443832462968
$ws.Range("L20").Value = 1.050077224734076
$ws.Range("M20").Value = 1.054825901166784
$ws.Range("N20").Value = 1.034860326753549

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.026225165204154
$ws.Range("D21").Value = 1.030400827492092
$ws.Range("E21").Value = 1.045378154674415
$ws.Range("F21").Value = 1.050079760340388
$ws.Range("I21").Value = 1.033471069185433
$ws.Range("J21").Value = 1.032882128247721
$ws.Range("K21").Value = 1.034001928909694
$ws.Range("L21").Value = 1.048923737137479
$ws.Range("M21").Value = 1.053608275382568
$ws.Range("N21").Value = 1.034348939343068

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.025652541410212
$ws.Range("D22").Value = 1.029991085985247
$ws.Range("E22").Value = 1.044522404105975
$ws.Range("F22").Value = 1.049184269921994
$ws.Range("I22").Value = 1.033356658392681
$ws.Range("J22").Value = 1.032560644406314
$ws.Range("K22").Value = 1.033723578429391
$ws.Range("L22").Value = 1.048198883739997
$ws.Range("M22").Value = 1.052843144041215
$ws.Range("N22").Value = 1.034026998957732

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.025956070486864
$ws.Range("D23").Value = 1.03020827319312
$ws.Range("E23").Value = 1.044975960313371
$ws.Range("F23").Value = 1.049658889599028
$ws.Range("I23").Value = 1.033417415886045
$ws.Range("J23").Value = 1.032731107688504
$ws.Range("K23").Value = 1.03387118516341
$ws.Range("L23").Value = 1.048583105532901
$ws.Range("M23").Value = 1.053248713332906
$ws.Range("N23").Value = 1.034197704317335

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.027151686577954
$ws.Range("D24").Value = 1.031063854565462
$ws.Range("E24").Value = 1.046763645016974
$ws.Range("F24").Value = 1.051529574021377
$ws.Range("I24").Value = 1.033654255328052
$ws.Range("J24").Value = 1.033401344764795
$ws.Range("K24").Value = 1.034451232355566
$ws.Range("L24").Value = 1.050096569778342
$ws.Range("M24").Value = 1.05484632230924
$ws.Range("N24").Value = 1.034868893207152

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.028540435975814
$ws.Range("D25").Value = 1.032057765981371
$ws.Range("E25").Value = 1.048842488889939
$ws.Range("F25").Value = 1.053704888786202
$ws.Range("I25").Value = 1.033924193980773
$ws.Range("J25").Value = 1.034177303954264
$ws.Range("K25").Value = 1.035122102673543
$ws.Range("L25").Value = 1.051854585677245
$ws.Range("M25").Value = 1.056702191205472
$ws.Range("N25").Value = 1.035645954347674
